$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix casing of the Vector3 data class name (was "vector3")
$ws.Range("F3").Value = "Vector3"

# Add new column H: field name "TestValue" of type "Vector3[]"
$ws.Range("H2").Value = "TestValue"
$ws.Range("H3").Value = "Vector3[]"

# Update selection to reflect the newly-edited cell
$ws.Range("H4").Select()
